$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 214 (Región Metropolitana,
# 2022-07-27). Every existing record from the old row 214 down to the old
# row 228 shifts down by one row, so the new last row becomes row 229
# (a duplicate of the former row 228). Shift bottom-up so we never
# overwrite a source row before it has been copied.
for ($r = 228; $r -ge 214; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Now populate the new record in row 214.
$ws.Range("D214").Value = 44769
$ws.Range("J214").Value = 500
$ws.Range("M214").Value = 10540
$ws.Range("O214").Value = "Región Metropolitana"
$ws.Range("P214").Value = 211
